# Update the 'want-to-go count' values in column F (index column).
# Sheet 1 = Exhibitions, Sheet 2 = Performances, Sheet 4 = All types (combined).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 5384
$ws.Range("F4").Value = 182
$ws.Range("F5").Value = 224
$ws.Range("F6").Value = 191
$ws.Range("F7").Value = 9013
$ws.Range("F8").Value = 280
$ws.Range("F9").Value = 656
$ws.Range("F10").Value = 19
$ws.Range("F11").Value = 2659
$ws.Range("F12").Value = 6356
$ws.Range("F13").Value = 2368
$ws.Range("F15").Value = 21
$ws.Range("F16").Value = 32
$ws.Range("F17").Value = 2559
$ws.Range("F20").Value = 6714
$ws.Range("F21").Value = 239
$ws.Range("F23").Value = 171
$ws.Range("F26").Value = 7451
$ws.Range("F29").Value = 250
$ws.Range("F30").Value = 49
$ws.Range("F32").Value = 14
$ws.Range("F34").Value = 40
$ws.Range("F36").Value = 24
$ws.Range("F37").Value = 42
$ws.Range("F38").Value = 61
$ws.Range("F39").Value = 2565
$ws.Range("F42").Value = 24
$ws.Range("F43").Value = 1141
$ws.Range("F45").Value = 580
$ws.Range("F46").Value = 3587
$ws.Range("F47").Value = 119
$ws.Range("F49").Value = 67

$ws = $wb.Worksheets.Item(2)
$ws.Range("F5").Value = 229
$ws.Range("F8").Value = 22

$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 5384
$ws.Range("F3").Value = 182
$ws.Range("F4").Value = 224
$ws.Range("F5").Value = 191
$ws.Range("F6").Value = 9013
$ws.Range("F7").Value = 280
$ws.Range("F8").Value = 656
$ws.Range("F9").Value = 19
$ws.Range("F11").Value = 2661
$ws.Range("F13").Value = 229
$ws.Range("F14").Value = 6356
$ws.Range("F15").Value = 2368
$ws.Range("F17").Value = 21
$ws.Range("F18").Value = 32
$ws.Range("F19").Value = 2559
$ws.Range("F21").Value = 22
$ws.Range("F24").Value = 6714
$ws.Range("F25").Value = 239
$ws.Range("F27").Value = 171
$ws.Range("F29").Value = 7451
$ws.Range("F31").Value = 250
$ws.Range("F32").Value = 49
$ws.Range("F33").Value = 14
$ws.Range("F35").Value = 24
$ws.Range("F36").Value = 42
$ws.Range("F38").Value = 61
$ws.Range("F39").Value = 2565
$ws.Range("F40").Value = 24
$ws.Range("F41").Value = 1141
$ws.Range("F43").Value = 580
$ws.Range("F45").Value = 3587
$ws.Range("F46").Value = 119
$ws.Range("F49").Value = 67

